# Redirect to the home page if quiz slug is invalid.
# Updates the "Report" sheet with refreshed (sanitized) quiz history rows,
# including new email values, updated correct/incorrect answer counts,
# and additional rows reflecting further quiz activity.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: History / Sanjib Roy / sanjibroy0098@gmail.com -- counts change
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 2

# Row 3: email updated, counts change
$ws.Cells.Item(3, 3).Value = "sanjibroy0099@gmail.com"
$ws.Cells.Item(3, 4).Value = 2
$ws.Cells.Item(3, 5).Value = 0

# Row 4: email updated, counts change
$ws.Cells.Item(4, 3).Value = "sanjibroy01@gmail.com"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1

# Row 5: email updated, counts change
$ws.Cells.Item(5, 3).Value = "sanjibroy@gmail.com"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1

# Row 6: quiz name corrected back to History, email updated, counts change
$ws.Cells.Item(6, 1).Value = "History"
$ws.Cells.Item(6, 2).Value = "Sanjib Roy"
$ws.Cells.Item(6, 3).Value = "sanjibroy0095@gmail.com"
$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1

# Row 7 (new): further History attempt
$ws.Cells.Item(7, 1).Value = "History"
$ws.Cells.Item(7, 2).Value = "Sanjib Roy"
$ws.Cells.Item(7, 3).Value = "sanjibronjjknkj@gmail.com"
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 2

# Row 8 (new): further History attempt
$ws.Cells.Item(8, 1).Value = "History"
$ws.Cells.Item(8, 2).Value = "Sanjib Roy"
$ws.Cells.Item(8, 3).Value = "sanjibroybjhbjhbhjbhj0098@gmail.com"
$ws.Cells.Item(8, 4).Value = 2
$ws.Cells.Item(8, 5).Value = 0

# Row 9 (new): a "Test Quiz" attempt
$ws.Cells.Item(9, 1).Value = "Test Quiz"
$ws.Cells.Item(9, 2).Value = "Sanjib Roy"
$ws.Cells.Item(9, 3).Value = "sanjibroy0098@gmail.com"
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 1

# Column C widened to fit the longest new email address.
$ws.Columns.Item(3).ColumnWidth = 41.14285714285714
